$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Objetivos:" row (row 10): previously held docente text by mistake,
# now holds the real course-objectives paragraph ---
$ws.Range("B10").Value = "Propiciar uma integração entre os elementos de estruturação da cidade, das variáveis ambientais e da malha urbana."
$ws.Range("C10").Value = "Propiciar uma integração entre os elementos de estruturação da cidade, das variáveis ambientais e da malha urbana."

# --- Insert a new row at 13 (pushes old rows 13-21 down to 14-22) and give
# it the "Docentes responsáveis:" value, which was previously missing ---
$ws.Rows.Item(13).Insert()

# Copy the formatting of an existing B/C pair down onto the new row 13 so the
# cell styles (wrap text etc.) match the rest of the table, then set values.
$ws.Range("B10:C10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B13").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C13").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

# --- Row 14 (was row 13) "Programa resumido:" used to say "Semestral";
# now holds the real short-syllabus text ---
$ws.Range("B14").Value = "Variável Ecológicano Ambiente Urbano; Enfoque Encômico e Impactos Ambientais."
$ws.Range("C14").Value = "Variável Ecológicano Ambiente Urbano; Enfoque Encômico e Impactos Ambientais."

# --- Row 16 (was row 15) "Programa:" used to hold a stray date value;
# now holds the full syllabus paragraph ---
$ws.Range("B16").Value = "Elementos para estruturação ambiental da cidade. Variável ecológica no ambiente das atividades urbanas. A questão ambiental no urbanismo. A questão ambiental sob o enfoque econômico. Noções de higiene e saúde ambiental. A urbanização e os impactos ocasionados, principal enfoque da drenagem urbana."
$ws.Range("C16").Value = "Elementos para estruturação ambiental da cidade. Variável ecológica no ambiente das atividades urbanas. A questão ambiental no urbanismo. A questão ambiental sob o enfoque econômico. Noções de higiene e saúde ambiental. A urbanização e os impactos ocasionados, principal enfoque da drenagem urbana."

# --- Row 19 (was row 18) "Método:" used to repeat the docente text;
# now holds the real teaching-method text ---
$ws.Range("B19").Value = "Aula expositiva e exercícios dirigidos."
$ws.Range("C19").Value = "Aula expositiva e exercícios dirigidos."

# --- Row 20 (was row 19) "Critério:" ---
$ws.Range("B20").Value = "Média ponderada de exercícios e provas."
$ws.Range("C20").Value = "Média ponderada de exercícios e provas."

# --- Row 21 (was row 20) "Norma de recuperação:" ---
$ws.Range("B21").Value = "Prova única com nota igual ou superior a 5,0."
$ws.Range("C21").Value = "Prova única com nota igual ou superior a 5,0."

# --- Row 22 (was row 21) "Bibliografia:" used to hold the recovery-rule
# text; now holds the real bibliography, and a brand-new row 22 is added
# with the actual bibliography paragraph ---
$bib = "valle, C.R. Qualidade ambiental: o desafio de ser competitivo protegendo o meio ambiente. Pioneira. 1995.`nDonaire, D.. Gestão ambiental na empresa. Atlas. 2a. edição. 1999.`nWinter, G.. Gestão e ambiente. Modelo prático de integração empresarial. Texto Editora, Lisboa. 1992.`nTucci, C.E., Porto, R.M., L.L. e Barros, M.T. org.. Drenagem Urbana. Ed. da Universidade e ABRH. 1995."
$ws.Range("B22").Value = $bib
$ws.Range("C22").Value = $bib
